$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matches")

# Missing home-team / away-team values for rows 50-64 (matches 49-64)
$teams = @{
    50 = @("PAL", "BOT")
    51 = @("SLB", "CHE")
    52 = @("PSG", "MIA")
    53 = @("FLA", "BAY")
    54 = @("INT", "FLU")
    55 = @("MCI", "HIL")
    56 = @("RMA", "JUV")
    57 = @("BVB", "CFM")
    58 = @("FLU", "HIL")
    59 = @("PAL", "CHE")
    60 = @("PSG", "BAY")
    61 = @("RMA", "BVB")
    62 = @("FLU", "CHE")
    63 = @("PSG", "RMA")
    64 = @("CHE", "PSG")
}

foreach ($row in $teams.Keys) {
    $pair = $teams[$row]
    $ws.Cells.Item($row, 8).Value = $pair[0]
    $ws.Cells.Item($row, 9).Value = $pair[1]
}

$wb.Save()
